$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pomiar rw")

# --- Row 3: quantity bumped 1 -> 2 (E3 = D3*4 recalculates automatically) ---
$ws.Range("D3").Value = 2

# --- Row 7: quantity bumped 12 -> 13 (E7 shared formula recalculates automatically) ---
$ws.Range("D7").Value = 13

# --- New font/style: red text used for the "przejrzeć masy i flagi" note (row 16) ---
# Creating this style before the hyperlink work below keeps it landing on the same
# cellXf slot used for B16 in the target workbook.
$ws.Range("B16").Value = "przejrzeć masy i flagi "
$ws.Range("B16").Font.Color = 255

# --- New component rows (11-15) ---
$ws.Range("B11").Value = "p mos"
$ws.Range("C11").Value = "AO4425"
$ws.Range("D11").Value = 2
$ws.Range("E11").Formula = "=D11*4"
$ws.Range("F11").Value = "AO4425"
$ws.Range("G11").Value = 8
$ws.Range("I11").Value = "SO8"

$ws.Range("B12").Value = "r 100k "
$ws.Range("D12").Value = 1
$ws.Range("E12").Formula = "=D12*4"

$ws.Range("B13").Value = "c 1n"

$ws.Range("B14").Value = "c 100n"

$ws.Range("B15").Value = "złącze bateri XC90?"

# --- J10 gains the ACS712 datasheet hyperlink (same visual style as J2/J3) ---
$ws.Hyperlinks.Add($ws.Range("J10"), "https://www.tme.eu/Document/75af175f80c090e9b8f9078a0b0b2409/ACS712.PDF")
$ws.Range("J2").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column B widened to fit the new text ---
$ws.Columns(2).ColumnWidth = 15.66

# --- Selection moved to C18 ---
$ws.Range("C18").Select()
